$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Belief Theory (BT) begins from a single primitive and models all change as the resolution of contradictions within a system.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Belief Theory (BT) begins from a single primitive and models change as arising from the resolution of contradictions within a system.",
    2)

$d.Content.Find.Execute(
    "Resolution is a sequence of atomic updates that eliminates contradiction.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Resolution is the elimination of a contradiction through a sequence of atomic updates.",
    2)

$d.Content.Find.Execute(
    "Information cost equals the updates required to eliminate the contradictions it introduces.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Information cost equals the number of atomic updates required to resolve the contradictions it introduces.",
    2)

$d.Content.Find.Execute(
    "Universal evolution follows resolution paths with minimal propagation distance.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Universal evolution proceeds along resolution paths that minimize propagation distance.",
    2)
